$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "Stimuli/vibraphonePlaceholder1.wav"
$ws.Range("B2").Value = "Stimuli/trigger_vibraphonePlaceholder1.wav"

$ws.Range("A8").Value = "Stimuli/vibraphonePlaceholder2.wav"
$ws.Range("B8").Value = "Stimuli/trigger_vibraphonePlaceholder2.wav"
